# Insert a new price-record row at row 35 (shifting all subsequent rows
# down by one, through the former last row 161 which becomes row 162),
# then populate the newly inserted row 35 with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("35:35").Insert()

$ws.Range("A35").Value = 10
$ws.Range("B35").Value = "Vega Modelo de Temuco"
$ws.Range("C35").Value = "La Araucanía"
$ws.Range("D35").Value = 44487
$ws.Range("E35").Value = 9
$ws.Range("F35").Value = 100112043
$ws.Range("G35").Value = "Pepino dulce"
$ws.Range("H35").Value = "Cultivar IV Región"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 80
$ws.Range("K35").Value = 19000
$ws.Range("L35").Value = 20000
$ws.Range("M35").Value = 19312
$ws.Range("N35").Value = "$/bandeja 18 kilos"
$ws.Range("O35").Value = "Provincia de Limarí"
$ws.Range("P35").Value = 1073
$ws.Range("Q35").Value = 18
$ws.Range("R35").Value = "Hortaliza"
